# Saldo_guide.xlsx update
# - Rename the data-extract sheet to reflect the new export timestamp
#   (20240925-085931 -> 20240926-085323)
# - Bump every "G" column date serial from 45560 (2024-09-25) to
#   45561 (2024-09-26) for all data rows (2..274)
# - Refresh the row-108 balances (D/E/H) with the newly reported figures

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet to match the new extraction run id
$wb.Worksheets.Item(1).Name = "IClientBalance-20240926-085323-"

# Update the date column (G) for every data row: 45560 -> 45561
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45561
}

# Row 108 carries revised totals in addition to the date bump
$ws.Range("D108").Value = 121993.15
$ws.Range("E108").Value = 44251.58
$ws.Range("H108").Value = 166244.73
